$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 201.45
$ws.Range("I33").Value = 198.93333
$ws.Range("J33").Value = 209
$ws.Range("K33").Value = 198.93333
$ws.Range("L33").Value = 209
$ws.Range("M33").Value = 30.06666999999999
$ws.Range("N33").Value = -667
$ws.Range("H53").Value = 230.64706
$ws.Range("J53").Value = 346.4
$ws.Range("L53").Value = 346.4
$ws.Range("N53").Value = -1620.4
$ws.Range("H106").Value = 3332.0557
$ws.Range("I106").Value = 3215.2666
$ws.Range("K106").Value = 3215.2666
$ws.Range("M106").Value = -2584.2666
$ws.Range("H112").Value = 3084.3
$ws.Range("J112").Value = 3112.4473
$ws.Range("L112").Value = 9337.341899999999
$ws.Range("N112").Value = -11553.3419
$ws.Range("H116").Value = 14102829
$ws.Range("J116").Value = 5156.1
$ws.Range("L116").Value = 5156.1
$ws.Range("N116").Value = -12040.1
$ws.Range("H141").Value = 3702.9285
$ws.Range("I141").Value = 2564.2
$ws.Range("K141").Value = 7692.599999999999
$ws.Range("M141").Value = -2512.599999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1035.7142
$ws.Range("J4").Value = 1025
$ws.Range("L4").Value = 1025
$ws.Range("N4").Value = -1257
$ws.Range("H32").Value = 6271.27
$ws.Range("I32").Value = 4625.043
$ws.Range("K32").Value = 4625.043
$ws.Range("M32").Value = -4338.043
$ws.Range("H45").Value = 3075.9656
$ws.Range("I45").Value = 2735.6316
$ws.Range("K45").Value = 2735.6316
$ws.Range("M45").Value = -2358.6316
$ws.Range("H56").Value = 8833.333000000001
$ws.Range("I56").Value = 3500
$ws.Range("J56").Value = 19500
$ws.Range("K56").Value = 3500
$ws.Range("L56").Value = 19500
$ws.Range("M56").Value = -2758
$ws.Range("N56").Value = -20984
$ws.Range("H74").Value = 11365250
$ws.Range("I74").Value = 14707201
$ws.Range("K74").Value = 14707201
$ws.Range("M74").Value = -14706327
$ws.Range("H77").Value = 11365250
$ws.Range("I77").Value = 14707201
$ws.Range("K77").Value = 73536005
$ws.Range("M77").Value = -73531637
$ws.Range("H97").Value = 517.4761999999999
$ws.Range("I97").Value = 531.8823
$ws.Range("K97").Value = 531.8823
$ws.Range("M97").Value = -35.88229999999999
$ws.Range("H122").Value = 2663.7812
$ws.Range("I122").Value = 2198.037
$ws.Range("K122").Value = 6594.110999999999
$ws.Range("M122").Value = -4144.110999999999
$ws.Range("H132").Value = 11825.058
$ws.Range("I132").Value = 14452.2
$ws.Range("K132").Value = 43356.60000000001
$ws.Range("M132").Value = -40826.60000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2852.6592
$ws.Range("I86").Value = 2575.7144
$ws.Range("J86").Value = 3337.3125
$ws.Range("K86").Value = 2575.7144
$ws.Range("L86").Value = 3337.3125
$ws.Range("M86").Value = -1452.7144
$ws.Range("N86").Value = -5583.3125
$ws.Range("H89").Value = 2852.6592
$ws.Range("I89").Value = 2575.7144
$ws.Range("J89").Value = 3337.3125
$ws.Range("K89").Value = 12878.572
$ws.Range("L89").Value = 16686.5625
$ws.Range("M89").Value = -7262.572
$ws.Range("N89").Value = -27918.5625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 6961.364
$ws.Range("I86").Value = 6938.2
$ws.Range("J86").Value = 6980.6665
$ws.Range("K86").Value = 6938.2
$ws.Range("L86").Value = 6980.6665
$ws.Range("M86").Value = -5815.2
$ws.Range("N86").Value = -9226.666499999999
$ws.Range("H89").Value = 6961.364
$ws.Range("I89").Value = 6938.2
$ws.Range("J89").Value = 6980.6665
$ws.Range("K89").Value = 34691
$ws.Range("L89").Value = 34903.3325
$ws.Range("M89").Value = -29075
$ws.Range("N89").Value = -46135.3325
$ws.Range("H97").Value = 5818.1816
$ws.Range("J97").Value = 6100
$ws.Range("L97").Value = 6100
$ws.Range("N97").Value = -8082
$ws.Range("H99").Value = 7236.3687
$ws.Range("I99").Value = 7375.905
$ws.Range("J99").Value = 7064
$ws.Range("K99").Value = 7375.905
$ws.Range("L99").Value = 7064
$ws.Range("M99").Value = -5877.905
$ws.Range("N99").Value = -10060
$ws.Range("H126").Value = 7236.3687
$ws.Range("I126").Value = 7375.905
$ws.Range("J126").Value = 7064
$ws.Range("K126").Value = 22127.715
$ws.Range("L126").Value = 21192
$ws.Range("M126").Value = -19657.715
$ws.Range("N126").Value = -26132
$ws.Range("H132").Value = 30305198
$ws.Range("I132").Value = 34189890
$ws.Range("J132").Value = 4601.6
$ws.Range("K132").Value = 102569670
$ws.Range("L132").Value = 13804.8
$ws.Range("M132").Value = -102567140
$ws.Range("N132").Value = -18864.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 942.5769
$ws.Range("I5").Value = 435.73334
$ws.Range("J5").Value = 1633.7273
$ws.Range("K5").Value = 1307.20002
$ws.Range("L5").Value = 4901.1819
$ws.Range("M5").Value = -1195.20002
$ws.Range("N5").Value = -5125.1819
$ws.Range("H68").Value = 1233.125
$ws.Range("I68").Value = 977.5
$ws.Range("K68").Value = 2932.5
$ws.Range("M68").Value = -2121.5
$ws.Range("H71").Value = 1233.125
$ws.Range("I71").Value = 977.5
$ws.Range("K71").Value = 8797.5
$ws.Range("M71").Value = -4741.5
$ws.Range("H135").Value = 942.5769
$ws.Range("I135").Value = 435.73334
$ws.Range("J135").Value = 1633.7273
$ws.Range("K135").Value = 3921.60006
$ws.Range("L135").Value = 14703.5457
$ws.Range("M135").Value = -1386.60006
$ws.Range("N135").Value = -19773.5457

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5941.7334
$ws.Range("I70").Value = 5755.75
$ws.Range("K70").Value = 5755.75
$ws.Range("M70").Value = -5485.75
$ws.Range("H73").Value = 5941.7334
$ws.Range("I73").Value = 5755.75
$ws.Range("K73").Value = 5755.75
$ws.Range("M73").Value = -4819.75
$ws.Range("H97").Value = 1673.5454
$ws.Range("I97").Value = 1507.4
$ws.Range("K97").Value = 1507.4
$ws.Range("M97").Value = -1011.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 17246418
$ws.Range("J40").Value = 18524602
$ws.Range("L40").Value = 18524602
$ws.Range("N40").Value = -18524874
$ws.Range("H61").Value = 3217.5833
$ws.Range("I61").Value = 3314
$ws.Range("K61").Value = 3314
$ws.Range("M61").Value = -3112
$ws.Range("H113").Value = 3217.5833
$ws.Range("I113").Value = 3314
$ws.Range("K113").Value = 3314
$ws.Range("M113").Value = -1144
$ws.Range("H122").Value = 7850.6577
$ws.Range("I122").Value = 4070.9443
$ws.Range("K122").Value = 12212.8329
$ws.Range("M122").Value = -9762.832900000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 501.55554
$ws.Range("I113").Value = 397.9565
$ws.Range("K113").Value = 1193.8695
$ws.Range("M113").Value = 976.1305
$ws.Range("H122").Value = 4550.647
$ws.Range("I122").Value = 3867.5881
$ws.Range("K122").Value = 11602.7643
$ws.Range("M122").Value = -9152.764299999999
$ws.Range("H132").Value = 8774594
$ws.Range("J132").Value = 3855.7334
$ws.Range("L132").Value = 11567.2002
$ws.Range("N132").Value = -16627.2002
$ws.Range("H136").Value = 4211.315
$ws.Range("I136").Value = 2938.158
$ws.Range("J136").Value = 7235.0625
$ws.Range("K136").Value = 8814.474
$ws.Range("L136").Value = 21705.1875
$ws.Range("M136").Value = -6264.474
$ws.Range("N136").Value = -26805.1875

